$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Analista de testes'
$ws.Range("B2").Value = 'Formel D Group'
$ws.Range("C2").Value = 'Hortolândia, São Paulo, Brasil (Presencial)'
$ws.Range("D2").Value = 'https://www.linkedin.com/jobs/view/4316648699/?eBP=CwEAAAGaA23jKkNVExKLwjxHAGBjA50PWwQJm7Fx7DdiqH5mJxphN1TKt0y7jv0VYFvbyMBbvJqsZtfUeyvh8d0ph8RXLVSQcqdr3x999G6qgGUtH9LP2vJQlBTqN7g37g_MV0VhCppQCToYRu1kIgYZSE0XbdCJ3ldH8EpmdU4kym7dyEe21jvE325Fy2VCxL-3GJ2-DYEHjKAsrS09ONHDpNKIFZvrxsuwnaVv18dUzYgav22RAY6VjDDarDKxaNQ8BbIfXxjLE_wA1Yorw2cfeQtESNDhsAj89ux6_j1Mih0PXqEv-aQvHPqdPmtw0nXOaRtE2OsGKcDDnW7OMe1AFxjK_6HSGQmaO9_csyUDtW1cHk0BK-lTtQtIE6HdvXiAXSSJFSEZssFdxcd_fcjnQ03e1sSyp21WrCCf1-9jTc_yPkJw9jzfSaENFQkbVkt9EqfSqnUxNofSVp_4aTo8FMF18GV2E9JS3TiszLgwpo1hUpDlXB_zpLBo24fOd6EU6yjmZaiRjA&refId=lbnt2cNJFf5kM%2FxHsC7d9A%3D%3D&trackingId=mLqNvWqSBSK%2Bh32dgP%2B6Pg%3D%3D&trk=flagship3_search_srp_jobs'

$ws.Range("A3").Value = 'Analista de Software Embarcado'
$ws.Range("B3").Value = 'AFRY'
$ws.Range("C3").Value = 'Curitiba, Paraná, Brasil (Presencial)'
$ws.Range("D3").Value = 'https://www.linkedin.com/jobs/view/4292388244/?eBP=CwEAAAGaA23jKorieug2KzePpeTJ7QXoTad1BS_QJ19ZUFkjZPY7OWIf9c4SPEYwowLkRc6jxa_qmbcCVLJBDtzEnt_vjTVHyYu_A8ycyYaUWmZb4rsPrOJxrtjE2M4MwlA7QyqtMhRK5dGPBTG9veazgdrGbzlFK0tO02wNFBLEH_43arEoyoHnfb2IVC7poVirYtU6MA3swkUc0w8ZJ5UtBGWCWCs5QP5VwU5dwN98_q9hD51txfyda4gBu1epLS937zaNjRmLJ87aaPYY_UEehAYVe5QAzP0QRnrvkglxwaH2McacNa20P_Ac7iCrO15kbq3qnZnsxo2-pisvPsVdo2JDsUG3FS1Y5e-IqkBwI_RA0Imd-rNbOPfPEiW05X6ZwT4_rZlqESj55MamLAJe64TkV6-6l9T3q4fjc4xfUrjcFzGodM0qU-UZ9jPghVpSvx06yfbZyEmVVQnwCJ-Q5Hqwf4O8HG5eu6ynbBl8_LzhJXX0OngaI4GRUVF0fhT17S-Fnmo9tS9c&refId=lbnt2cNJFf5kM%2FxHsC7d9A%3D%3D&trackingId=qX89QF2WypEYi9xm3e%2Fr5w%3D%3D&trk=flagship3_search_srp_jobs'

$ws.Range("A4").Value = 'Especialista I Em Desenvolvimento De Software'
$ws.Range("B4").Value = 'TOTVS'
$ws.Range("C4").Value = 'Joinville, Santa Catarina, Brasil (Remoto)'
$ws.Range("D4").Value = 'https://www.linkedin.com/jobs/view/4314916144/?eBP=CwEAAAGaA23jKrrhFD2-P3tQpLlx15pMo-3V0UkGnxIfI4jeX_sNd4LKZlLTu68v3eNAC93DpnZ0I3yaENP2azOhIM2CkHLPfs4pv18nre7XOZBQQFKtttG-SarwqidlUvb624P6ncO9FPDKNtEGMvPQperDwmYXYCPwEcWqzLJ3L9yYFjbbvdp3XQ0moKNqD4pGEu7qbfnetBFXic2-KR3XM7LqwplMNSE8uotnaoCZGgVjsvaVbgqUv4eS8BUjnTQcRob7BadRqwaj_aKdc-YptXQl6KNKq3XxluF8aLf56q1BtzJ8TsYZ6psUk1SdHXAvDUiSre8zGZjPto_DNLHmz29DI4hjCHQs1mI3b_womIojIqc0y71M7lmK8IggXOF4JhKsipFIv4n4Sa4sX74GeLn7IzDm3rgfGFWxe6UVqME9kLC-RjLKNQxeePOSotnveS0S20WMU1IlEt6yHXjmJCxUMOb5Ri2ceU8u61bbeLR7rtN-vUncCde0Nzfx-SMlrTcPRjYlZA&refId=lbnt2cNJFf5kM%2FxHsC7d9A%3D%3D&trackingId=%2Fcn6MbQWo9Ui2%2BVsQ7Hkig%3D%3D&trk=flagship3_search_srp_jobs'

$ws.Range("A5").Value = 'ANALISTA DESENVOLVIMENTO HARDWARE (PD40E7)'
$ws.Range("B5").Value = 'Marelli'
$ws.Range("C5").Value = 'Hortolândia, São Paulo, Brasil (Híbrido)'
$ws.Range("D5").Value = 'https://www.linkedin.com/jobs/view/4317357720/?eBP=BUDGET_EXHAUSTED_JOB&refId=lbnt2cNJFf5kM%2FxHsC7d9A%3D%3D&trackingId=d9wtTWEXo9%2FOk1daq33qAQ%3D%3D&trk=flagship3_search_srp_jobs'

$ws.Range("A6").Value = 'DESENVOLVEDOR DE HARDWARE AUTOMOMOTIVO PLENO | HÍBRIDO (Curitiba)'
$ws.Range("B6").Value = 'Mobi7 Localiza'
$ws.Range("C6").Value = 'Curitiba, Paraná, Brasil (Híbrido)'
$ws.Range("D6").Value = 'https://www.linkedin.com/jobs/view/4314908154/?eBP=CwEAAAGaA23jKn2vKEFYv3eCW9b3haFSH-J_pM2FHgZ9uWG6UIER-R3u-VEjK3-B1mMn2mjwwY0TAqFs8ZQcR7W1LHY8kTY3KH9Gsrl5G_LH_S4A-sdSRxC65itJWPbYewr78riFtEl2blOXqiQuBeWDZaqwTGdUmrNSTMcxU5p0xfy8lMi5eAsmZXz0d9D4g2AXkId9FjVkMuFscfHl1MZNGQGUFtPFtEH6ECrASI8HjBVW7ZPrjVM2l5iIj_HgqPmqSpHJKenvSPqvToan9gZXyek7qID9RtDzpRsfIbCZayeNcjUkc2EvfArxD679QqJ4cyPEEYMK6XTn5BXfvTT7OMSBAqcrJMq9wRLl6wRg5_NDWXNUTiT2WsLHT5Qo05vl9sJwIawzkBalADdiNPt_qWW93z_-GNsDS0B_ygHgIxiZcyWsDFvKOaNUqc8mAyVbcKcNT8X_GbQsUAGS4OXnQe4_O8FEZ80qTOeakEU&refId=lbnt2cNJFf5kM%2FxHsC7d9A%3D%3D&trackingId=4LZNlHcTmZEWGLTfYj1F0A%3D%3D&trk=flagship3_search_srp_jobs'

$ws.Range("A7").Value = 'ANALISTA P&D JR'
$ws.Range("B7").Value = 'Honda Brasil'
$ws.Range("C7").Value = 'Sumaré, São Paulo, Brasil (Presencial)'
$ws.Range("D7").Value = 'https://www.linkedin.com/jobs/view/4298033073/?eBP=NOT_ELIGIBLE_FOR_CHARGING&refId=lbnt2cNJFf5kM%2FxHsC7d9A%3D%3D&trackingId=Evv0cGJiO1rylxDcAzbe6Q%3D%3D&trk=flagship3_search_srp_jobs'

$ws.Range("A8").Value = 'JUNIOR EMBEDDED ENGINEER - AGRICULTURAL AUTOMATION (32478)'
$ws.Range("B8").Value = 'Bosch Brasil'
$ws.Range("C8").Value = 'Curitiba e Região (Presencial)'
$ws.Range("D8").Value = 'https://www.linkedin.com/jobs/view/4308286750/?eBP=CwEAAAGaA23jKoHE5TFBkZ-uuDwb60pVomW96MrOjRKYnnBLh6ckli9Tw64g1SUhiP_lLCUct003NbsZCqYyC7T8a0f5FIPF8bEyiakKknGW91EmCp4wSYcI7Bj89TStZCudbdIwcmcULotrQotwmh01BcAju8ACsLZq4-BCceeiaP43w9RagAHvUub5NF9VJUWwh-KenDVVh22SvsWJC1YXSWKSyp3qfhxiqU428dgoec_7gLs_f4jzjW9tQKrEE7eDGZUX4UQ68xocKj2espE52ZMx1Gx3SN5NllObo4gPxFwwIZqFAz7OWBaOVQWPJeaDEOi_5GC6xCTeUSLPGUfErR--gYruPiB5p6niuIn2gTs5kZN9ZAYaaVV61KxAuuW_xNcOAmXe0zxt2gYV-E7Bqxpt87DwDmATTEhaS1GuzivVpxJR_o_iPiEdxgL3MyHqW9RbkQ1w0uo7j_uXZ3qsF22rXOIfpm3ueN5WNaFHVMOJT4CiZBrkS_7Tsr-F5Dx_F29aCn67b7oP&refId=lbnt2cNJFf5kM%2FxHsC7d9A%3D%3D&trackingId=IFgaQoKwqC1bK2bcl1sORQ%3D%3D&trk=flagship3_search_srp_jobs'
